$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "32÷9="

$cell = $t.Cell(1, 2)
$cell.Range.Text = "62÷7="

$cell = $t.Cell(1, 3)
$cell.Range.Text = "50÷9="

$cell = $t.Cell(1, 4)
$cell.Range.Text = "12÷7="

$cell = $t.Cell(1, 5)
$cell.Range.Text = "17÷7="

$cell = $t.Cell(5, 1)
$cell.Range.Text = "44÷6="

$cell = $t.Cell(5, 2)
$cell.Range.Text = "19÷7="

$cell = $t.Cell(5, 3)
$cell.Range.Text = "21÷9="

$cell = $t.Cell(5, 4)
$cell.Range.Text = "97÷3="

$cell = $t.Cell(5, 5)
$cell.Range.Text = "30÷4="

$cell = $t.Cell(9, 1)
$cell.Range.Text = "77÷2="

$cell = $t.Cell(9, 2)
$cell.Range.Text = "16÷6="

$cell = $t.Cell(9, 3)
$cell.Range.Text = "88÷5="

$cell = $t.Cell(9, 4)
$cell.Range.Text = "84÷4="

$cell = $t.Cell(9, 5)
$cell.Range.Text = "26÷2="

$cell = $t.Cell(13, 1)
$cell.Range.Text = "24÷2="

$cell = $t.Cell(13, 2)
$cell.Range.Text = "45÷7="

$cell = $t.Cell(13, 3)
$cell.Range.Text = "46÷7="

$cell = $t.Cell(13, 4)
$cell.Range.Text = "64÷6="

$cell = $t.Cell(13, 5)
$cell.Range.Text = "79÷6="

$cell = $t.Cell(17, 1)
$cell.Range.Text = "96÷2="

$cell = $t.Cell(17, 2)
$cell.Range.Text = "50÷8="

$cell = $t.Cell(17, 3)
$cell.Range.Text = "71÷4="

$cell = $t.Cell(17, 4)
$cell.Range.Text = "54÷7="

$cell = $t.Cell(17, 5)
$cell.Range.Text = "52÷3="
